$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-04-20"

# Update the header label for the current-month column (column B)
$ws.Range("B1").Value = "April 2022 (through April 20)"

# --- Neighborhood-by-month count updates for 2022-04-28 data refresh ---

# Englewood
$ws.Range("Z3").Value = 3

# Chatham
$ws.Range("B11").Value = 1

# Portage Park
$ws.Range("Z13").Value = 1

# Lincoln Park
$ws.Range("B24").Value = 3
$ws.Range("N24").Value = 2

# South Shore
$ws.Range("B26").Value = 3

# Roseland
$ws.Range("F32").Value = 4

# Gage Park
$ws.Range("F47").Value = 1
$ws.Range("R47").Value = 1

# Kenwood
$ws.Range("N53").Value = 1

# Archer Heights
$ws.Range("F58").Value = 1

# Ashburn
$ws.Range("R60").Value = 1

# Hegewisch
$ws.Range("J70").Value = 1

# Montclare
$ws.Range("F77").Value = 2

# Old Town
$ws.Range("E82").Value = 2

# United Center - remove erroneous prior count
$ws.Range("B91").ClearContents()

# West Pullman
$ws.Range("F93").Value = 2
